$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.998.87"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.234.72"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'245.66"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").Value = "'76.10"
$ws.Range("E7").Value = "  +5.77%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'40.86"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "2.573.80"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "'14.80"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "'0.856"
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").Value = "2.229.57"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "41.951.35"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "0.0₃0976"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "'6.10"
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("D21").Value = "'71.23"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'2.21"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'229.75"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'3.71"
$ws.Range("E25").Value = "  -5.37%  "
$ws.Range("D26").Value = "'11.11"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").Value = "'2.31"
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").Value = "'7.24"
$ws.Range("E28").Value = "  +12.90%  "
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "'169.24"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").Value = "'20.44"
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("D32").Value = "'0.0855"
$ws.Range("E32").Value = "  +4.84%  "
$ws.Range("D33").Value = "'33.19"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").Value = "'4.79"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").Value = "'13.27"
$ws.Range("E39").Value = "  -4.81%  "
$ws.Range("D40").Value = "'2.19"
$ws.Range("E40").Value = "  -5.93%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'113.20"
$ws.Range("E42").Value = "  +16.27%  "
$ws.Range("E43").Value = "  -6.50%  "
$ws.Range("D44").Value = "'59.85"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").Value = "'8.77"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'1.12"
$ws.Range("E48").Value = "  -4.56%  "
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "'4.21"
$ws.Range("E50").Value = "  -13.96%  "
$ws.Range("E51").Value = "  -1.31%  "
